$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows (2-5) down to (3-6)
$ws.Rows("2:2").Insert()

# Fill in the new row 2 with the Minciencias postdoc info
# (order chosen to reproduce the shared-string table ordering of the target file)
$ws.Range("C2").Value = "Ministerio de Ciencia Tecnología e Innovación - Minciencias"
$ws.Range("D2").Value = "Bogotá, Colombia"
$ws.Range("B2").Value = "2023 - 2025"
$ws.Range("A2").Value = "Convocatoria de Estancias Posdoctorales Orientadas por Misiones"
$ws.Range("E2").Value = "\textbf{Proyecto: } La necesidad de generar procesos de reparación social a las mujeres víctimas y sobrevivientes de violencias sexuales en el marco del conflicto armado desde el quehacer periodístico. Diversas propuestas de tratamiento según contextos"

$ws.Range("E2").Select()
